# Reorders the "section header" rows in the demographics table (Sheet1)
# down by one row within each labelled block (Education Completed,
# Household, Age Migrated, Migration Cohort, Acculturation), matching the
# layout already used by the "Age" block at the top of the sheet.
#
# Concretely each block [header, data1, data2, ...] is rotated up by one
# row so the header becomes the last row of its block:
#   rows  7- 8 : Education Completed block (2 rows)
#   rows 12-14 : Household block (3 rows)
#   rows 17-20 : Age Migrated block (4 rows)
#   rows 22-26 : Migration Cohort block (5 rows)
#   rows 27-32 : Acculturation block (6 rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final contents (after the re-order) for every row that actually changes.
# $null clears a cell (used for the rows that become bare section headers).
$rowData = @{
    7  = @('Less than Primary', '0.35', '0.14', '0.25', '0.1', '0.24', '0.08', '0.09', '0.07', '0.03', '0.01', '0.03')
    8  = @('Education Completed', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    12 = @('Household Size', '3.32', '2.31', '2.87', '2.46', '3.14', '2.73', '2.65', '2.37', '2.06', '1.95', '2.27')
    13 = @('Lives Alone', '0.16', '0.32', '0.24', '0.27', '0.18', '0.21', '0.21', '0.28', '0.39', '0.29', '0.31')
    14 = @('Household', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    17 = @('Less than 15', '0.11', '0.29', '0.06', '0.17', '0.06', '0.07', '0.13', '-', '-', '-', '-')
    18 = @('15 - 24', '0.31', '0.3', '0.19', '0.16', '0.23', '0.2', '0.2', '-', '-', '-', '-')
    19 = @('25 - 49', '0.45', '0.25', '0.52', '0.42', '0.56', '0.51', '0.47', '-', '-', '-', '-')
    20 = @('Age Migrated', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    22 = @('Before 1965', '0.15', '0.41', '0.08', '0.25', '0.09', '0.11', '0.21', '-', '-', '-', '-')
    23 = @('1965 - 1979', '0.45', '0.28', '0.33', '0.35', '0.37', '0.34', '0.31', '-', '-', '-', '-')
    24 = @('1980 - 1999', '0.33', '0.18', '0.43', '0.28', '0.46', '0.36', '0.35', '-', '-', '-', '-')
    25 = @('After 1999', '0.12', '0.14', '0.2', '0.21', '0.14', '0.23', '0.17', '-', '-', '-', '-')
    26 = @('Migration Cohort', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    27 = @('Citizen', '0.54', '-', '0.68', '0.8', '0.66', '0.71', '0.74', '-', '-', '-', '-')
    28 = @('English Speakers', '0.73', '0.91', '0.69', '0.76', '0.82', '0.87', '0.93', '0.99', '1', '1', '1')
    29 = @('Cognitive Difficulty', '0.09', '0.15', '0.13', '0.13', '0.08', '0.07', '0.08', '0.11', '0.12', '0.08', '0.11')
    30 = @('Independence Difficulty', '0.15', '0.21', '0.18', '0.18', '0.13', '0.1', '0.14', '0.15', '0.18', '0.12', '0.15')
    31 = @('N', '79658', '23021', '8980', '21242', '20064', '25003', '352960', '120724', '313063', '3165675', '94162')
    32 = @('Acculturation', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
}

foreach ($r in $rowData.Keys | Sort-Object) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($r, $col)
        $v = $vals[$i]
        if ($null -eq $v) {
            $cell.Value = $null
        } elseif ($col -eq 1) {
            # Column A values are plain labels - never numeric-looking.
            $cell.Value = $v
        } else {
            # Prefix with an apostrophe so values that look numeric
            # ("0.35", "79658", "-") are stored as text, matching the
            # original table's formatting, then drop the resulting
            # "quote prefix" style so the cell keeps the default style.
            $cell.Value = "'" + $v
            $cell.Style = "Normal"
        }
    }
}
